# Weekly roll-forward of the Perejil (parsley) price series.
# The tracked date block lives in rows 154:211 (pairs of "Primera"/"Segunda"
# quality rows per reporting date, newest date first). A new week's data is
# inserted at the top, every older row slides down two rows, and the two
# oldest rows that fall off the bottom of the previously-seen window are
# appended as brand-new rows 212:213 (growing the sheet by 2 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 154:211 down by two rows. Excel shifts the existing data to
# 156:213 automatically (and rows 212:213 end up holding exactly what used
# to be rows 210:211 - the oldest pair - which is exactly what the new
# workbook needs).
$ws.Rows("154:155").Insert()

# Seed the freshly-inserted (blank) rows 154:155 with the same row shape as
# the prior newest entry, which has now moved to 156:157.
$ws.Range("A156:R157").Copy()
$ws.Range("A154").PasteSpecial()

# Stamp the new reporting date on the newly-added pair.
$ws.Range("D154").Value = 45006
$ws.Range("D155").Value = 45006
